# Corrects issue in "Modelling" column.
# Fills column M ("Modelling") with "Yes" for rows 116 through 165,
# matching the value already present in column N for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Target Species")

for ($r = 116; $r -le 165; $r++) {
    $ws.Cells.Item($r, 13).Value = "Yes"   # Column M = 13
}

# Update the active selection to reflect where the edit took place.
$ws.Range("K155").Select()
